$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price/Volume columns so that
# numeric-looking strings (e.g. "1.004", "0.3822") are stored as text,
# matching the original inlineStr/text cell type instead of being
# auto-converted to floating point numbers by Excel.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '28.066.43'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '1.777.53'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").Value = '339.56'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("D7").Value = '0.3822'
$ws.Range("E7").Value = '  -2.46%  '
$ws.Range("D8").Value = '0.3416'
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").Value = '46.79'
$ws.Range("E9").Value = '  -3.06%  '
$ws.Range("E10").Value = '  -4.43%  '
$ws.Range("D11").Value = '0.07387'
$ws.Range("E11").Value = '  -1.49%  '
$ws.Range("D12").Value = '23.22'
$ws.Range("E12").Value = '  +5.80%  '
$ws.Range("D13").Value = '1.004'
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").Value = '6.377'
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").Value = '7.422'
$ws.Range("E15").Value = '  +3.62%  '
$ws.Range("D16").Value = '1.777.44'
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").Value = '0.00001075'
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("D18").Value = '0.06662'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '82.46'
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("D21").Value = '17.36'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = '6.401'
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("D23").Value = '28.076.68'
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("D24").Value = '12.07'
$ws.Range("E24").Value = '  -3.10%  '
$ws.Range("D25").Value = '2.379'
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("D26").Value = '1.446'
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("D27").Value = '20.71'
$ws.Range("E27").Value = '  -2.49%  '
$ws.Range("D28").Value = '2.410'
$ws.Range("E28").Value = '  -3.95%  '
$ws.Range("D29").Value = '154.12'
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").Value = '1.980.27'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").Value = '134.30'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").Value = '4.031'
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("D33").Value = '6.056'
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = '0.08904'
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("E35").Value = '  -2.87%  '
$ws.Range("D36").Value = '0.02399'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").Value = '0.6843'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").Value = '0.06392'
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("D39").Value = '5.291'
$ws.Range("E39").Value = '  -3.13%  '
$ws.Range("D40").Value = '0.2163'
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("D41").Value = '1.239'
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").Value = '1.499'
$ws.Range("E42").Value = '  -7.44%  '
$ws.Range("D43").Value = '8.222'
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("D44").Value = '14.24'
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("D46").Value = '0.6265'
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("D47").Value = '3.865'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").Value = '133.04'
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("D49").Value = '2.070'
$ws.Range("E49").Value = '  -3.17%  '
$ws.Range("D50").Value = '0.07505'
$ws.Range("E50").Value = '  +4.13%  '
$ws.Range("D51").Value = '1.208'
$ws.Range("E51").Value = '  +4.24%  '

# Restore the default (unstyled) cell style so the text-format change
# above does not leave a visible style index on the cells.
$rng.Style = "Normal"
